$wb = $excel.ActiveWorkbook

# --- RQ1 sheet: fill in the previously-missing Jabref (row 7) CF/Revision
#     counts for both Regular Clones and Micro Clones, which makes the
#     E7/H7 ratio formulas resolve instead of #DIV/0!, and ripples into the
#     Total row (row 8) sums/ratios. ---
$ws1 = $wb.Worksheets.Item("RQ1")
$ws1.Range("C7").Value = 165
$ws1.Range("D7").Value = 65
$ws1.Range("F7").Value = 316
$ws1.Range("G7").Value = 65

# Remove the leftover "Average" formula textbox/shape that used to sit over
# RQ1 (it referenced data that's now obsolete).
for ($i = $ws1.Shapes.Count; $i -ge 1; $i--) {
    $ws1.Shapes.Item($i).Delete()
}

# --- RQ4 sheet: drop the stray AVERAGE row (row 15) that only ever produced
#     #DIV/0! errors. ---
$ws4 = $wb.Worksheets.Item("RQ4")
$ws4.Rows("15:15").Delete()

# --- Update which sheet/cell is active & selected: RQ1 (selection E8)
#     becomes the active tab instead of RQ4 (selection moves to G15). ---
$ws4.Range("G15").Select() | Out-Null
$ws1.Activate()
$ws1.Range("E8").Select() | Out-Null
